$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("concept")

# Update the Weibull start values (C_start: gas / pv)
$ws.Range("D10").Value = 20
$ws.Range("E10").Value = 20

# Update the new-capacity inputs (C_new: gas / pv) for every year row 15-34
$ws.Range("G15:G34").Value = 2
$ws.Range("H15:H34").Value = 2

# Restore the last active selection on the sheet
$ws.Range("Q38").Select()
